# Auto-generated edit script applying numeric updates to Sheets/Belias_Profits.xlsx
# as described in the commit diff. Values were recomputed upstream (e.g. market
# price refreshes) and this script reproduces the same per-row profit figures.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1038.8
$ws.Range("I32").Value = 995
$ws.Range("J32").Value = 1049.75
$ws.Range("K32").Value = 995
$ws.Range("L32").Value = 1049.75
$ws.Range("M32").Value = -669
$ws.Range("N32").Value = -1701.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H10").Value = 15601.8
$ws.Range("I10").Value = 2252
$ws.Range("J10").Value = 24501.666
$ws.Range("K10").Value = 2252
$ws.Range("L10").Value = 24501.666
$ws.Range("M10").Value = -2082
$ws.Range("N10").Value = -24841.666

$ws.Range("H11").Value = 5000000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H12").Value = 850
$ws.Range("J12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("N12").Value = -1346

$ws.Range("H13").Value = 32000000
$ws.Range("I13").Value = 32000000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 32000000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -31999856
$ws.Range("N13").ClearContents()

$ws.Range("H14").Value = 950
$ws.Range("I14").Value = 950
$ws.Range("K14").Value = 950
$ws.Range("M14").Value = -775

$ws.Range("H16").Value = 1500
$ws.Range("J16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2574

$ws.Range("H19").Value = 12833.333
$ws.Range("J19").Value = 12833.333
$ws.Range("L19").Value = 12833.333
$ws.Range("N19").Value = -13291.333

$ws.Range("H27").Value = 10600
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 10600
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 10600
$ws.Range("N27").Value = -10968
$ws.Range("M27").ClearContents()

$ws.Range("H30").Value = 6820.9
$ws.Range("I30").Value = 2704.5
$ws.Range("J30").Value = 7850
$ws.Range("K30").Value = 2704.5
$ws.Range("L30").Value = 7850
$ws.Range("M30").Value = -2554.5
$ws.Range("N30").Value = -8150

$ws.Range("H32").Value = 3979.1807
$ws.Range("I32").Value = 2876.125
$ws.Range("J32").Value = 7694.737
$ws.Range("K32").Value = 2876.125
$ws.Range("L32").Value = 7694.737
$ws.Range("M32").Value = -2589.125
$ws.Range("N32").Value = -8268.737000000001

$ws.Range("H46").Value = 71285.2
$ws.Range("I46").Value = 5114
$ws.Range("J46").Value = 170542
$ws.Range("K46").Value = 5114
$ws.Range("L46").Value = 170542
$ws.Range("M46").Value = -4795
$ws.Range("N46").Value = -171180

$ws.Range("H74").Value = 57175.668
$ws.Range("I74").Value = 201386.8
$ws.Range("J74").Value = 1709.8462
$ws.Range("K74").Value = 201386.8
$ws.Range("L74").Value = 1709.8462
$ws.Range("M74").Value = -200512.8
$ws.Range("N74").Value = -3457.8462

$ws.Range("H77").Value = 57175.668
$ws.Range("I77").Value = 201386.8
$ws.Range("J77").Value = 1709.8462
$ws.Range("K77").Value = 1006934
$ws.Range("L77").Value = 8549.231
$ws.Range("M77").Value = -1002566
$ws.Range("N77").Value = -17285.231

$ws.Range("H122").Value = 1671.7037
$ws.Range("I122").Value = 1074.8667
$ws.Range("J122").Value = 2417.75
$ws.Range("K122").Value = 3224.6001
$ws.Range("L122").Value = 7253.25
$ws.Range("M122").Value = -774.6001000000001
$ws.Range("N122").Value = -12153.25

$ws.Range("H132").Value = 2598.7666
$ws.Range("I132").Value = 2616.0625
$ws.Range("J132").Value = 2579
$ws.Range("K132").Value = 7848.1875
$ws.Range("L132").Value = 7737
$ws.Range("M132").Value = -5318.1875
$ws.Range("N132").Value = -12797

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 6464
$ws.Range("I25").Value = 1206.4
$ws.Range("J25").Value = 9750
$ws.Range("K25").Value = 1206.4
$ws.Range("L25").Value = 9750
$ws.Range("M25").Value = -971.4000000000001
$ws.Range("N25").Value = -10220

$ws.Range("H54").Value = 6842.5557
$ws.Range("I54").Value = 694.3333
$ws.Range("K54").Value = 694.3333
$ws.Range("M54").Value = -210.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 60001580
$ws.Range("I31").Value = 100001380
$ws.Range("J31").Value = 33335050
$ws.Range("K31").Value = 100001380
$ws.Range("L31").Value = 33335050
$ws.Range("M31").Value = -100001085
$ws.Range("N31").Value = -33335640

$ws.Range("H34").Value = 60001580
$ws.Range("I34").Value = 100001380
$ws.Range("J34").Value = 33335050
$ws.Range("K34").Value = 100001380
$ws.Range("L34").Value = 33335050
$ws.Range("M34").Value = -100001178
$ws.Range("N34").Value = -33335454

$ws.Range("H62").Value = 2555.111
$ws.Range("I62").Value = 2247.6924
$ws.Range("J62").Value = 3354.4
$ws.Range("K62").Value = 2247.6924
$ws.Range("L62").Value = 3354.4
$ws.Range("M62").Value = -1623.6924
$ws.Range("N62").Value = -4602.4

$ws.Range("H65").Value = 2555.111
$ws.Range("I65").Value = 2247.6924
$ws.Range("J65").Value = 3354.4
$ws.Range("K65").Value = 11238.462
$ws.Range("L65").Value = 16772
$ws.Range("M65").Value = -8118.462
$ws.Range("N65").Value = -23012

$ws.Range("H134").Value = 2050
$ws.Range("I134").Value = 2130
$ws.Range("J134").Value = 1823.3334
$ws.Range("K134").Value = 6390
$ws.Range("L134").Value = 5470.0002
$ws.Range("M134").Value = -3855
$ws.Range("N134").Value = -10540.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 140
$ws.Range("I11").Value = 80
$ws.Range("J11").Value = 230
$ws.Range("K11").Value = 240
$ws.Range("L11").Value = 690
$ws.Range("M11").Value = -100
$ws.Range("N11").Value = -970

$ws.Range("H140").Value = 2879.96
$ws.Range("I140").Value = 785.6429000000001
$ws.Range("K140").Value = 2356.9287
$ws.Range("M140").Value = 2823.0713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 62.344826
$ws.Range("I2").Value = 51.842106
$ws.Range("J2").Value = 82.3
$ws.Range("K2").Value = 51.842106
$ws.Range("L2").Value = 82.3
$ws.Range("M2").Value = 61.157894
$ws.Range("N2").Value = -308.3

$ws.Range("H122").Value = 102650
$ws.Range("I122").Value = 144657.14
$ws.Range("J122").Value = 4633.3335
$ws.Range("K122").Value = 433971.42
$ws.Range("L122").Value = 13900.0005
$ws.Range("M122").Value = -431521.42
$ws.Range("N122").Value = -18800.0005

$ws.Range("H132").Value = 2665.6216
$ws.Range("I132").Value = 2423.4
$ws.Range("J132").Value = 3703.7144
$ws.Range("K132").Value = 7270.200000000001
$ws.Range("L132").Value = 11111.1432
$ws.Range("M132").Value = -4740.200000000001
$ws.Range("N132").Value = -16171.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3501.2727
$ws.Range("I7").Value = 3213
$ws.Range("J7").Value = 4270
$ws.Range("K7").Value = 3213
$ws.Range("L7").Value = 4270
$ws.Range("M7").Value = -3101
$ws.Range("N7").Value = -4494

$ws.Range("H34").Value = 6950
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 6950
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 6950
$ws.Range("N34").Value = -7294
$ws.Range("M34").ClearContents()

$ws.Range("H126").Value = 3501.2727
$ws.Range("I126").Value = 3213
$ws.Range("J126").Value = 4270
$ws.Range("K126").Value = 9639
$ws.Range("L126").Value = 12810
$ws.Range("M126").Value = -7169
$ws.Range("N126").Value = -17750

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 12161.4
$ws.Range("I61").Value = 6500
$ws.Range("J61").Value = 17822.8
$ws.Range("K61").Value = 6500
$ws.Range("L61").Value = 17822.8
$ws.Range("M61").Value = -6208
$ws.Range("N61").Value = -18406.8

$ws.Range("H122").Value = 6282211
$ws.Range("I122").Value = 9617362
$ws.Range("K122").Value = 28852086
$ws.Range("M122").Value = -28849636

$ws.Range("H126").Value = 333333860
$ws.Range("I126").Value = 250000540
$ws.Range("J126").Value = 500000500
$ws.Range("K126").Value = 750001620
$ws.Range("L126").Value = 1500001500
$ws.Range("M126").Value = -749999150
$ws.Range("N126").Value = -1500006440

$ws.Range("H135").Value = 45000
$ws.Range("J135").Value = 45000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -55140

$ws.Range("H136").Value = 83339490
$ws.Range("I136").Value = 111118264
$ws.Range("J136").Value = 3166.6667
$ws.Range("K136").Value = 333354792
$ws.Range("L136").Value = 9500.000100000001
$ws.Range("M136").Value = -333352242
$ws.Range("N136").Value = -14600.0001

$ws.Range("H137").Value = 300000
$ws.Range("J137").Value = 300000
$ws.Range("L137").Value = 300000
$ws.Range("N137").Value = -310200
